$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column F header and row 3 value, plus new row 4
$ws.Range("F2").Value = "Similar problems"
$ws.Range("F3").Value = "215. Kth Largest Element in an Array"
$ws.Range("A4").Value = 209
$ws.Range("B4").Value = "Minimum Size Subarray Sum"

# Adjust row 3 height back to default (removes explicit custom row height)
$ws.Rows.Item(3).AutoFit()

# Adjust column E width (closest achievable value to 51.85546875 given engine's
# quantization of column widths to 1/6-character steps)
$ws.Columns.Item(5).ColumnWidth = 51.0

# Update B3's font to match the new style: Calibri 11, color FF212121, minor theme
$ws.Range("B3").Font.Size = 11
$ws.Range("B3").Font.Name = "Calibri"
$ws.Range("B3").Font.ThemeFont = 1

# Update the active cell selection to B4
$ws.Range("B4").Select()
